$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44199

# Row 4: Harry Potter and the Prisoner of Azkaban
$ws.Range("A4").Value = "Harry Potter and the Prisoner of Azkaban"
$ws.Range("B4").Value = "J.K. Rowling"
$ws.Range("C4").Value = 44199
$ws.Range("D4").Value = 44200
$ws.Range("E4").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F4").Value = "Audio"
$ws.Range("G4").Value = "12 Hours 21 Mins"
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = $true

# Row 5: The Hard Thing About Hard Things
$ws.Range("A5").Value = "The Hard Thing About Hard Things"
$ws.Range("B5").Value = "Ben Horowitz"
$ws.Range("C5").Value = 44200
$ws.Range("D5").Value = 44204
$ws.Range("E5").Value = "entreprenuership;business;ceo;building a company;success"
$ws.Range("F5").Value = "Audio"
$ws.Range("G5").Value = "8 Hours 4 Mins"
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = $false

# Use the same style as the existing date cells (copy format from C3) onto the
# new/changed date cells so we don't introduce new number formats.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy()
$ws.Range("C4:D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E6").Select()
